$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68f52f3ff0ef848b7178db5cedc6664655d83d34/e2e/7e760aac-2fba-487a-83c8-0770a6daff2a.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e7ade6d32217a2dbff7f443c4a73b5214e589f3e/e2e/7e760aac-2fba-487a-83c8-0770a6daff2a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68f52f3ff0ef848b7178db5cedc6664655d83d34/e2e/7e760aac-2fba-487a-83c8-0770a6daff2a.md."

# ----- zh-cn sheet -----
$wsZh = $wb.Worksheets.Item("zh-cn")

# widen the "Error Detail" column (P / 16) to fit the long message
$wsZh.Columns.Item(16).ColumnWidth = 39.17

$wsZh.Range("I8").Value = "7e760aac-2fba-487a-83c8-0770a6daff2a.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $latestUrl, "", "", "7e760aac-2fba-487a-83c8-0770a6daff2a.md")
$wsZh.Range("J8").Value = "7e760aac-2fba-487a-83c8-0770a6daff2a.b681b75fe143645720764fec195184ce998003bd.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-08-17 22:42:28"
$wsZh.Range("P8").Value = $errorDetail

# ----- de-de sheet -----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.17

$wsDe.Range("I8").Value = "7e760aac-2fba-487a-83c8-0770a6daff2a.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $latestUrl, "", "", "7e760aac-2fba-487a-83c8-0770a6daff2a.md")
$wsDe.Range("J8").Value = "7e760aac-2fba-487a-83c8-0770a6daff2a.b681b75fe143645720764fec195184ce998003bd.de-de.xlf"
$wsDe.Range("K8").Value = "2016-08-17 22:42:35"
$wsDe.Range("P8").Value = $errorDetail
